$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = "r567"
$ws.Range("B10").Value = "trudy"
$ws.Range("C10").Value = "works as expected"
$ws.Range("D10").Value = "2025-09-30 13:44:09"
